$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove existing hyperlinks up-front; they will be re-created below
# pointing at the (shifted) rows so we don't end up with stale refs.
$ws.Hyperlinks.Delete()

# --- Insert a brand-new row 20 ("Liver" / Hepatocellular Carcinoma ...),
# pushing the former rows 20-34 down to 21-35.
$ws.Rows("20:20").Insert()

# --- Populate the freshly inserted row 20.
$ws.Range("A20").Value = "Liver"

# B20 needs the special "s=2" formatting that B7 already carries, so copy
# B7's formatting over first and then overwrite the text.
$ws.Range("B7").Copy($ws.Range("B20"))
$ws.Range("B20").Value = "Hepatocellular Carcinoma with Intra-lesional Air (Superinfected HCC)"

$ws.Range("C20").Value = "Clip 1 B-mode"

# D20 looks like a hyperlink cell (style s=1, "Collegamento ipertestuale")
# but, per the source diff, it is NOT wired up as an actual hyperlink -
# copy the style from another D-column hyperlink-styled cell, then just
# set the text (no Hyperlinks.Add for this one).
$ws.Range("D19").Copy($ws.Range("D20"))
$ws.Range("D20").Value = "https://youtu.be/4V95TgIMrbE"

# --- Re-create the hyperlinks for every D-column cell below the new row,
# using the same target URLs as before the insert (Excel's row-insert
# does not repoint hyperlink refs itself, so we redo them from scratch).
# Hyperlinks.Add stamps its own "applied" hyperlink style variant, so we
# immediately paste back just the formats from an existing hyperlink cell
# (D19) to keep every linked D cell on the workbook's original style index.
function Add-Link($cellRef, $target) {
    $cell = $ws.Range($cellRef)
    $ws.Hyperlinks.Add($cell, $target) | Out-Null
    $ws.Range("D19").Copy()
    $cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

Add-Link "D4"  "https://youtu.be/zxTC0YBY2RY "
Add-Link "D27" "https://youtu.be/xBfd04F4Ni8 "
Add-Link "D11" "https://youtu.be/91M82AIMyu0 "
Add-Link "D32" "https://youtu.be/qushjTAy6XQ "
Add-Link "D29" "https://youtu.be/pc-vbxSRTbs "
Add-Link "D21" "https://youtu.be/DjI1kEnzfSQ "
Add-Link "D28" "https://youtu.be/JvwODCASLYQ "
Add-Link "D22" "https://youtu.be/U3ydTsRwxok "
Add-Link "D14" "https://youtu.be/15o_Km86IzM "
Add-Link "D33" "https://youtu.be/_FckFwJwynI "
Add-Link "D30" "https://youtu.be/Axbee4vjNtU"
Add-Link "D16" "https://youtu.be/RhSUFLTmTl4"
Add-Link "D7"  "https://youtu.be/2kRZcpi70Aw "
Add-Link "D34" "https://youtu.be/z_oaRVxRz5s "
Add-Link "D5"  "https://youtu.be/K2Wbg7BgXy4 "
Add-Link "D3"  "https://youtu.be/ZXwd0gwHEkQ "
Add-Link "D31" "https://youtu.be/VJdnjrAAO-4"
Add-Link "D2"  "https://youtu.be/kdZO1IPuOIw"
Add-Link "D35" "https://youtu.be/S45odD2wQOQ"

# --- Restore the selection to the newly added row's D cell, matching the
# author's final selection in the workbook.
$ws.Range("D20").Select()
